# Scheduled market-data refresh: update computed price/profit columns (H:N)
# for the rows whose item market data changed, across all crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 32: Automata for the People / Crab Oil
$ws.Cells.Item(32, 8).Value = 734
$ws.Cells.Item(32, 9).Value = 450
$ws.Cells.Item(32, 10).Value = 1302
$ws.Cells.Item(32, 11).Value = 450
$ws.Cells.Item(32, 12).Value = 1302
$ws.Cells.Item(32, 13).Value = -124
$ws.Cells.Item(32, 14).Value = -1954

# ALC row 64: Forged from the Void / Void Glue
$ws.Cells.Item(64, 8).Value = 3575.5334
$ws.Cells.Item(64, 9).Value = 3503.8462
$ws.Cells.Item(64, 10).Value = 3673.6316
$ws.Cells.Item(64, 11).Value = 3503.8462
$ws.Cells.Item(64, 12).Value = 3673.6316
$ws.Cells.Item(64, 13).Value = -3255.8462
$ws.Cells.Item(64, 14).Value = -4169.631600000001

# ALC row 67: Dodging the Draft (L) / Void Glue
$ws.Cells.Item(67, 8).Value = 3575.5334
$ws.Cells.Item(67, 9).Value = 3503.8462
$ws.Cells.Item(67, 10).Value = 3673.6316
$ws.Cells.Item(67, 11).Value = 3503.8462
$ws.Cells.Item(67, 12).Value = 3673.6316
$ws.Cells.Item(67, 13).Value = -2645.8462
$ws.Cells.Item(67, 14).Value = -5389.631600000001

# ALC row 116: Growing Up / Growth Formula Kappa
$ws.Cells.Item(116, 8).Value = 5921660
$ws.Cells.Item(116, 9).Value = 15392155
$ws.Cells.Item(116, 10).Value = 2600
$ws.Cells.Item(116, 11).Value = 15392155
$ws.Cells.Item(116, 12).Value = 2600
$ws.Cells.Item(116, 13).Value = -15388713
$ws.Cells.Item(116, 14).Value = -9484

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 36741.484
$ws.Cells.Item(137, 9).Value = 1735.5883
$ws.Cells.Item(137, 10).Value = 86333.164
$ws.Cells.Item(137, 11).Value = 5206.7649
$ws.Cells.Item(137, 12).Value = 258999.492
$ws.Cells.Item(137, 13).Value = -2656.7649
$ws.Cells.Item(137, 14).Value = -264099.492

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32: Ingot We Trust / Steel Ingot
$ws.Cells.Item(32, 8).Value = 16399386
$ws.Cells.Item(32, 9).Value = 21743002
$ws.Cells.Item(32, 10).Value = 12294.066
$ws.Cells.Item(32, 11).Value = 21743002
$ws.Cells.Item(32, 12).Value = 12294.066
$ws.Cells.Item(32, 13).Value = -21742715
$ws.Cells.Item(32, 14).Value = -12868.066

# ARM row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 1844.25
$ws.Cells.Item(61, 9).Value = 1406.25
$ws.Cells.Item(61, 10).Value = 2939.25
$ws.Cells.Item(61, 11).Value = 1406.25
$ws.Cells.Item(61, 12).Value = 2939.25
$ws.Cells.Item(61, 13).Value = -1194.25
$ws.Cells.Item(61, 14).Value = -3363.25

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 2485.1052
$ws.Cells.Item(132, 9).Value = 2432.1667
$ws.Cells.Item(132, 10).Value = 2575.8572
$ws.Cells.Item(132, 11).Value = 7296.500100000001
$ws.Cells.Item(132, 12).Value = 7727.571599999999
$ws.Cells.Item(132, 13).Value = -4766.500100000001
$ws.Cells.Item(132, 14).Value = -12787.5716

# ARM row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 1844.25
$ws.Cells.Item(136, 9).Value = 1406.25
$ws.Cells.Item(136, 10).Value = 2939.25
$ws.Cells.Item(136, 11).Value = 4218.75
$ws.Cells.Item(136, 12).Value = 8817.75
$ws.Cells.Item(136, 13).Value = -1668.75
$ws.Cells.Item(136, 14).Value = -13917.75

# ARM row 138: Don't Ask about the Rivets / Titanium Gold Helm of Casting
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 22: Riveting Run / Iron Rivets
$ws.Cells.Item(22, 8).Value = 127.333336
$ws.Cells.Item(22, 9).Value = 127.333336
$ws.Cells.Item(22, 11).Value = 127.333336
$ws.Cells.Item(22, 13).Value = 45.666664

# BSM row 107: The Gold Experience / Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 1809.5483
$ws.Cells.Item(107, 9).Value = 2102.35
$ws.Cells.Item(107, 10).Value = 1277.1818
$ws.Cells.Item(107, 11).Value = 2102.35
$ws.Cells.Item(107, 12).Value = 1277.1818
$ws.Cells.Item(107, 13).Value = -182.3499999999999
$ws.Cells.Item(107, 14).Value = -5117.1818

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 2013.6666
$ws.Cells.Item(134, 9).Value = 1926
$ws.Cells.Item(134, 10).Value = 2151.4285
$ws.Cells.Item(134, 11).Value = 5778
$ws.Cells.Item(134, 12).Value = 6454.2855
$ws.Cells.Item(134, 13).Value = -3243
$ws.Cells.Item(134, 14).Value = -11524.2855

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Cells.Item(58, 8).Value = 1239.1384
$ws.Cells.Item(58, 9).Value = 755.38
$ws.Cells.Item(58, 10).Value = 2851.6667
$ws.Cells.Item(58, 11).Value = 755.38
$ws.Cells.Item(58, 12).Value = 2851.6667
$ws.Cells.Item(58, 13).Value = -552.38
$ws.Cells.Item(58, 14).Value = -3257.6667

# CRP row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Cells.Item(134, 8).Value = 1470.9683
$ws.Cells.Item(134, 9).Value = 929.32074
$ws.Cells.Item(134, 10).Value = 4341.7
$ws.Cells.Item(134, 11).Value = 2787.96222
$ws.Cells.Item(134, 12).Value = 13025.1
$ws.Cells.Item(134, 13).Value = -252.9622199999999
$ws.Cells.Item(134, 14).Value = -18095.1

# CRP row 136: Turali Quality / Dark Mahogany Lumber
$ws.Cells.Item(136, 8).Value = 1239.1384
$ws.Cells.Item(136, 9).Value = 755.38
$ws.Cells.Item(136, 10).Value = 2851.6667
$ws.Cells.Item(136, 11).Value = 2266.14
$ws.Cells.Item(136, 12).Value = 8555.000100000001
$ws.Cells.Item(136, 13).Value = 283.8600000000001
$ws.Cells.Item(136, 14).Value = -13655.0001

$ws = $wb.Worksheets.Item("CUL")
# CUL row 2: Pork Is a Salty Food / Table Salt
$ws.Cells.Item(2, 8).Value = 40.25
$ws.Cells.Item(2, 9).Value = 10
$ws.Cells.Item(2, 10).Value = 55.375
$ws.Cells.Item(2, 11).Value = 60
$ws.Cells.Item(2, 12).Value = 332.25
$ws.Cells.Item(2, 13).Value = 53
$ws.Cells.Item(2, 14).Value = -558.25

# CUL row 5: What a Sap / Maple Syrup
$ws.Cells.Item(5, 8).Value = 498.73685
$ws.Cells.Item(5, 9).Value = 491.67648
$ws.Cells.Item(5, 10).Value = 558.75
$ws.Cells.Item(5, 11).Value = 1475.02944
$ws.Cells.Item(5, 12).Value = 1676.25
$ws.Cells.Item(5, 13).Value = -1363.02944
$ws.Cells.Item(5, 14).Value = -1900.25

# CUL row 33: Cooking with Gas / Chicken Stock
$ws.Cells.Item(33, 8).Value = 123.30769
$ws.Cells.Item(33, 9).Value = 113.75
$ws.Cells.Item(33, 10).Value = 138.6
$ws.Cells.Item(33, 11).Value = 682.5
$ws.Cells.Item(33, 12).Value = 831.5999999999999
$ws.Cells.Item(33, 13).Value = -399.5
$ws.Cells.Item(33, 14).Value = -1397.6

# CUL row 92: Oh No Udon / Gyr Abanian Flour
$ws.Cells.Item(92, 8).Value = 800
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 13).ClearContents()

# CUL row 127: A Stickler for Carrots / Carrot Nibbles
$ws.Cells.Item(127, 8).Value = 650
$ws.Cells.Item(127, 10).Value = 650
$ws.Cells.Item(127, 12).Value = 1950
$ws.Cells.Item(127, 14).Value = -11870

# CUL row 132: More Mezcal / Cooking Mezcal
$ws.Cells.Item(132, 8).Value = 10399543
$ws.Cells.Item(132, 9).Value = 847
$ws.Cells.Item(132, 10).Value = 14732333
$ws.Cells.Item(132, 11).Value = 7623
$ws.Cells.Item(132, 12).Value = 132590997
$ws.Cells.Item(132, 13).Value = -5093
$ws.Cells.Item(132, 14).Value = -132596057

# CUL row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Cells.Item(135, 8).Value = 498.73685
$ws.Cells.Item(135, 9).Value = 491.67648
$ws.Cells.Item(135, 10).Value = 558.75
$ws.Cells.Item(135, 11).Value = 4425.08832
$ws.Cells.Item(135, 12).Value = 5028.75
$ws.Cells.Item(135, 13).Value = -1890.08832
$ws.Cells.Item(135, 14).Value = -10098.75

$ws = $wb.Worksheets.Item("GSM")
# GSM row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Cells.Item(126, 8).Value = 3848134.2
$ws.Cells.Item(126, 9).Value = 5883792.5
$ws.Cells.Item(126, 10).Value = 3002.4443
$ws.Cells.Item(126, 11).Value = 17651377.5
$ws.Cells.Item(126, 12).Value = 9007.332900000001
$ws.Cells.Item(126, 13).Value = -17648907.5
$ws.Cells.Item(126, 14).Value = -13947.3329

# GSM row 132: On Board for Lar / Lar Ingot
$ws.Cells.Item(132, 8).Value = 4589.8613
$ws.Cells.Item(132, 9).Value = 4793.143
$ws.Cells.Item(132, 10).Value = 3878.375
$ws.Cells.Item(132, 11).Value = 14379.429
$ws.Cells.Item(132, 12).Value = 11635.125
$ws.Cells.Item(132, 13).Value = -11849.429
$ws.Cells.Item(132, 14).Value = -16695.125

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Cells.Item(132, 8).Value = 1556.13
$ws.Cells.Item(132, 9).Value = 1685.7294
$ws.Cells.Item(132, 10).Value = 821.73334
$ws.Cells.Item(132, 11).Value = 5057.1882
$ws.Cells.Item(132, 12).Value = 2465.20002
$ws.Cells.Item(132, 13).Value = -2527.1882
$ws.Cells.Item(132, 14).Value = -7525.20002

# LTW row 136: Respect for Br'aax / Br'aax Leather
$ws.Cells.Item(136, 8).Value = 1646.238
$ws.Cells.Item(136, 9).Value = 1280.619
$ws.Cells.Item(136, 10).Value = 2743.0952
$ws.Cells.Item(136, 11).Value = 3841.857
$ws.Cells.Item(136, 12).Value = 8229.285600000001
$ws.Cells.Item(136, 13).Value = -1291.857
$ws.Cells.Item(136, 14).Value = -13329.2856

$ws = $wb.Worksheets.Item("WVR")
# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 1549.4468
$ws.Cells.Item(132, 9).Value = 1019.7727
$ws.Cells.Item(132, 10).Value = 2015.56
$ws.Cells.Item(132, 11).Value = 3059.3181
$ws.Cells.Item(132, 12).Value = 6046.68
$ws.Cells.Item(132, 13).Value = -529.3181
$ws.Cells.Item(132, 14).Value = -11106.68
